$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 200
$ws.Range("B3").Value = 650
$ws.Range("B4").Value = 600
$ws.Range("B5").Value = 200
$ws.Range("B6").Value = 91
$ws.Range("B7").Value = 500
$ws.Range("B8").Value = 175
